# Word COM-interop script implementing the target diff:
#  1. First paragraph ("Hello, world") gains an en-US <w:lang> on its
#     paragraph mark rPr and on both runs' rPr.
#  2. The trailing empty paragraph is replaced by a 3x2 table
#     (style "TableGrid") followed by a new empty paragraph that also
#     carries the en-US <w:lang> on its paragraph mark rPr.
#  3. A "TableGrid" table style (built-in "Table Grid") is added to
#     styles.xml.
#
# Word doesn't expose a "patch raw XML of an arbitrary part" primitive,
# but Range.InsertXML() called on the *whole document* Range re-applies
# a full flat-OPC package (every <pkg:part/>, not just document.xml) in
# one shot, which is exactly the lever we need to touch styles.xml too.

$d = $word.ActiveDocument

$full = $d.WordOpenXML

# ---------------------------------------------------------------------
# 1) "Hello, world" paragraph: add en-US language everywhere.
# ---------------------------------------------------------------------

$old1 = '<w:p w14:paraId="37996B07" w14:textId="77777777" w:rsidR="00B66BAF" w:rsidRDefault="00B66BAF" w:rsidP="00B66BAF"><w:r><w:t xml:space="preserve">Hello, </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>world</w:t></w:r></w:p>'

$new1 = '<w:p w14:paraId="37996B07" w14:textId="77777777" w:rsidR="00B66BAF" w:rsidRDefault="00B66BAF" w:rsidP="00B66BAF"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Hello, </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>world</w:t></w:r></w:p>'

if ($full.IndexOf($old1) -lt 0) {
    throw "anchor #1 (Hello-world paragraph) not found"
}
$full = $full.Replace($old1, $new1)

# ---------------------------------------------------------------------
# 2) Replace the trailing empty paragraph with a table + empty para.
# ---------------------------------------------------------------------

$oldEmptyPara = '<w:p w14:paraId="54EFA04D" w14:textId="77777777" w:rsidR="00E1117F" w:rsidRDefault="00E1117F"/>'

$tableXml = '<w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2831"/><w:gridCol w:w="2831"/><w:gridCol w:w="2832"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="2831" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Abc</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2831" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2832" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2831" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2831" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>123</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2832" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The answer is: 42</w:t></w:r></w:p></w:tc></w:tr></w:tbl>'

$newTrailingPara = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'

if ($full.IndexOf($oldEmptyPara) -lt 0) {
    throw "anchor #2 (trailing empty paragraph) not found"
}
$full = $full.Replace($oldEmptyPara, ($tableXml + $newTrailingPara))

# ---------------------------------------------------------------------
# 3) Add the "TableGrid" table style to styles.xml.
# ---------------------------------------------------------------------

$styleXml = '<w:style w:type="table" w:styleId="TableGrid"><w:name w:val="Table Grid"/><w:basedOn w:val="TableNormal"/><w:uiPriority w:val="39"/><w:rsid w:val="00063ECB"/><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:tblPr><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders></w:tblPr></w:style>'

$stylesClose = '</w:styles>'
$closeIdx = $full.IndexOf($stylesClose)
if ($closeIdx -lt 0) {
    throw "anchor #3 (</w:styles>) not found"
}
$full = $full.Substring(0, $closeIdx) + $styleXml + $full.Substring($closeIdx)

# ---------------------------------------------------------------------
# Apply the whole rewritten flat-OPC package back onto the document.
# ---------------------------------------------------------------------

$d.Content.InsertXML($full)

Write-Output "applied"
